# Weekly update: insert a new price record for Choclo at
# "Vega Monumental Concepción" before the existing row 43, shifting the
# remaining historical rows (old 43-58) down by one row (new 44-59).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 43; existing rows 43:58 shift to 44:59.
$ws.Rows.Item(43).Insert()

# Populate the newly inserted row 43 with the new weekly record.
$ws.Cells.Item(43, 1).Value = 11
$ws.Cells.Item(43, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(43, 3).Value = "Bíobío"
$ws.Cells.Item(43, 4).Value = 44523
$ws.Cells.Item(43, 5).Value = 8
$ws.Cells.Item(43, 6).Value = 100112024
$ws.Cells.Item(43, 7).Value = "Choclo"
$ws.Cells.Item(43, 8).Value = "Dulce o Americano"
$ws.Cells.Item(43, 9).Value = "Primera"
$ws.Cells.Item(43, 10).Value = 100
$ws.Cells.Item(43, 11).Value = 24000
$ws.Cells.Item(43, 12).Value = 25000
$ws.Cells.Item(43, 13).Value = 24500
$ws.Cells.Item(43, 14).Value = "`$/malla 70 unidades"
$ws.Cells.Item(43, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(43, 16).Value = 350
$ws.Cells.Item(43, 17).Value = 70
$ws.Cells.Item(43, 18).Value = "Hortaliza"
